# Update the answers in the two-digit x two-digit multiplication table.
# The underlying diff removes the first answer cell of row 1 (shifting the
# remaining answers left) and appends one new answer cell at the end of
# row 1; every other cell keeps its position but gets a new answer value.
# Because every <w:tc> in the table shares identical formatting
# (tcPr/pPr/rPr), that structural change is equivalent to simply writing
# the new value into each of row 1's five cell positions -- so the whole
# edit can be done as straightforward per-cell text replacement.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "82×44=3608"
$t.Cell(1, 2).Range.Text = "54×90=4860"
$t.Cell(1, 3).Range.Text = "69×68=4692"
$t.Cell(1, 4).Range.Text = "63×55=3465"
$t.Cell(1, 5).Range.Text = "71×79=5609"

$t.Cell(5, 1).Range.Text = "90×57=5130"
$t.Cell(5, 2).Range.Text = "56×66=3696"
$t.Cell(5, 3).Range.Text = "37×18=666"
$t.Cell(5, 4).Range.Text = "78×29=2262"
$t.Cell(5, 5).Range.Text = "99×53=5247"

$t.Cell(10, 1).Range.Text = "91×69=6279"
$t.Cell(10, 2).Range.Text = "71×39=2769"
$t.Cell(10, 3).Range.Text = "86×49=4214"
$t.Cell(10, 4).Range.Text = "70×58=4060"
$t.Cell(10, 5).Range.Text = "75×76=5700"

$t.Cell(15, 1).Range.Text = "31×87=2697"
$t.Cell(15, 2).Range.Text = "51×47=2397"
$t.Cell(15, 3).Range.Text = "45×67=3015"
$t.Cell(15, 4).Range.Text = "19×64=1216"
$t.Cell(15, 5).Range.Text = "49×60=2940"

$t.Cell(20, 1).Range.Text = "76×14=1064"
$t.Cell(20, 2).Range.Text = "35×39=1365"
$t.Cell(20, 3).Range.Text = "60×97=5820"
$t.Cell(20, 4).Range.Text = "28×75=2100"
$t.Cell(20, 5).Range.Text = "25×39=975"
